# The commit swaps the deck's visual theme (ppt/theme/theme1.xml, the
# theme used by the one-and-only slide master) from the colourful
# "Integral" / "Red Violet" palette over to the plain default
# "Office Theme" palette. (ppt/theme/theme2.xml, which only the notes
# master points at, is not reachable through the PowerPoint object
# model and is left alone.)
#
# PowerPoint exposes per-slot theme colour editing through
# ThemeColorScheme.Colors(i).RGB -- that is the supported way to
# recolour a theme in place (renaming the theme / colour-scheme is not
# exposed by the object model, so only the RGB values are changed).

function HexToRgb($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in PowerPoint's ThemeColorScheme.Colors()
# slot order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgb $officeTheme[$i - 1]
}
